# "first running version for single 10be nuclide"
# On the "Composition and other" sheet, insert a new column before the
# existing "soil mass (g/cm2)" column (column C) to hold a new
# "CDF_err" field, with a sample value of 0.05 for the Test row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C; existing C/D/E (soil mass, scaling model,
# method) shift right to D/E/F. Excel carries over the formatting of
# the column to the left (B), matching the header/data styles used
# elsewhere in the row.
$ws.Columns.Item(3).Insert()

# New header + value for the inserted column.
$ws.Cells.Item(1, 3).Value = "CDF_err"
$ws.Cells.Item(2, 3).Value = 0.05

# Move the active selection to the new column, as in the edited file.
$ws.Range("C3").Select()
